$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the two data rows (74a19e4f.. and e43b8c80..) swap places.
# Row 2 becomes the e43b8c80 file (still "Handed back" status).
# Row 3 becomes the 74a19e4f file, now "Ready for handoff".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md"
$wsOverview.Range("A3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 04:53:29"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md", "", "", "e2e\e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md", "", "", "e2e\74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row swap, plus the handoff/handback file + datetime
# columns are updated, and row 3 (74a19e4f, now "Ready for handoff") gets a
# stale-handback error message.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md"
$wsZh.Range("G2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.08399a3d2c5f9ce10db2d70c76b41bb9e98e258d.zh-cn.xlf"
$wsZh.Range("I2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md"
$wsZh.Range("J2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.08399a3d2c5f9ce10db2d70c76b41bb9e98e258d.zh-cn.xlf"

$wsZh.Range("A3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.f710d0af1e39c37892691ebba7083365344ef8b9.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-04 04:53:24"
$wsZh.Range("I3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md"
$wsZh.Range("J3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.f710d0af1e39c37892691ebba7083365344ef8b9.zh-cn.xlf"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7bc042ed5e7a88b1d7a8f836651dbfdecef3dafd/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md."

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md", "", "", "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/98a778b6b98195b1a652f88a5ecc4a1cde995c71/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md", "", "", "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md", "", "", "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/98a778b6b98195b1a652f88a5ecc4a1cde995c71/e2e/e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md", "", "", "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md")

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# Sheet "de-de": identical pattern to zh-cn, but with de-de.xlf file names
# and its own handoff/handback datetimes.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md"
$wsDe.Range("G2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.08399a3d2c5f9ce10db2d70c76b41bb9e98e258d.de-de.xlf"
$wsDe.Range("I2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md"
$wsDe.Range("J2").Value = "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.08399a3d2c5f9ce10db2d70c76b41bb9e98e258d.de-de.xlf"

$wsDe.Range("A3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.f710d0af1e39c37892691ebba7083365344ef8b9.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-04 04:53:29"
$wsDe.Range("I3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md"
$wsDe.Range("J3").Value = "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.f710d0af1e39c37892691ebba7083365344ef8b9.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7bc042ed5e7a88b1d7a8f836651dbfdecef3dafd/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md."

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md", "", "", "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/585a4ad3fb62cd7903be9d708cc5b8c1531ebcab/e2e/74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md", "", "", "e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96f00d7f5438eb258169e9dadb2952da56be30a0/e2e/e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md", "", "", "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/585a4ad3fb62cd7903be9d708cc5b8c1531ebcab/e2e/e43b8c80-924c-4502-8377-1a1cd7bcc6f6.md", "", "", "74a19e4f-7c42-4b7e-adde-6ad79c96eb65.md")

$wsDe.Columns.Item(16).ColumnWidth = 39.17
